$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new Excel date serial value for column A (shift from 1st of quarter-start month
# to the 15th of the following month, per the "Quarterly indexing esoteric bug-fix operation" commit).
$dateUpdates = @{
    2 = 32370
    3 = 32462
    4 = 32554
    5 = 32643
    6 = 32735
    7 = 32827
    8 = 32919
    9 = 33008
    10 = 33100
    11 = 33192
    12 = 33284
    13 = 33373
    14 = 33465
    15 = 33557
    16 = 33649
    17 = 33739
    18 = 33831
    19 = 33923
    20 = 34015
    21 = 34104
    22 = 34196
    23 = 34288
    24 = 34380
    25 = 34469
    26 = 34561
    27 = 34653
    28 = 34745
    29 = 34834
    30 = 34926
    31 = 35018
    32 = 35110
    33 = 35200
    34 = 35292
    35 = 35384
    36 = 35476
    37 = 35565
    38 = 35657
    39 = 35749
    40 = 35841
    41 = 35930
    42 = 36022
    43 = 36114
    44 = 36206
    45 = 36295
    46 = 36387
    47 = 36479
    48 = 36571
    49 = 36661
    50 = 36753
    51 = 36845
    52 = 36937
    53 = 37026
    54 = 37118
    55 = 37210
    56 = 37302
    57 = 37391
    58 = 37483
    59 = 37575
    60 = 37667
    61 = 37756
    62 = 37848
    63 = 37940
    64 = 38032
    65 = 38122
    66 = 38214
    67 = 38306
    68 = 38398
    69 = 38487
    70 = 38579
    71 = 38671
    72 = 38763
    73 = 38852
    74 = 38944
    75 = 39036
    76 = 39128
    77 = 39217
    78 = 39309
    79 = 39401
    80 = 39493
    81 = 39583
    82 = 39675
    83 = 39767
    84 = 39859
    85 = 39948
    86 = 40040
    87 = 40132
    88 = 40224
    89 = 40313
    90 = 40405
    91 = 40497
    92 = 40589
    93 = 40678
    94 = 40770
    95 = 40862
    96 = 40954
    97 = 41044
    98 = 41136
    99 = 41228
    100 = 41320
    101 = 41409
    102 = 41501
    103 = 41593
    104 = 41685
    105 = 41774
    106 = 41866
    107 = 41958
    108 = 42050
    109 = 42139
    110 = 42231
    111 = 42323
    112 = 42415
    113 = 42505
    114 = 42597
    115 = 42689
    116 = 42781
    117 = 42870
    118 = 42962
    119 = 43054
    120 = 43146
    121 = 43235
    122 = 43327
    123 = 43419
    124 = 43511
    125 = 43600
    126 = 43692
    127 = 43784
    128 = 43876
    129 = 43966
    130 = 44058
    131 = 44150
    132 = 44242
    133 = 44331
    134 = 44423
    135 = 44515
    136 = 44607
    137 = 44696
    138 = 44788
    139 = 44880
    140 = 44972
    141 = 45061
    142 = 45153
    143 = 45245
    144 = 45337
    145 = 45427
    146 = 45519
    147 = 45611
    148 = 45703
    149 = 45792
    150 = 45884
}

foreach ($row in $dateUpdates.Keys) {
    $ws.Cells.Item($row, 1).Value = $dateUpdates[$row]
}

